# Simulador module improvements: update header labels (add units / currency
# hints), change the reference date sample value, widen the data columns to
# fit the new (longer) headers, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: keep each header in its original column, just reword it ---
$ws.Range("A1").Value = "DATA DE REFERÊNCIA"
$ws.Range("B1").Value = "CONSUMO (m³)"
$ws.Range("C1").Value = "REUTILIZADO (m³)"
$ws.Range("D1").Value = "(R$) VALOR PAGO"

# --- Sample reference date: 01/01/2000 -> 01/01/1999 ---
$ws.Range("A2").Value = 36161

# --- Widen columns B:D so the new, longer headers fit ---
$ws.Columns.Item(2).ColumnWidth = 17.833333333333332
$ws.Columns.Item(3).ColumnWidth = 19.0
$ws.Columns.Item(4).ColumnWidth = 18.5

# --- Move the active selection ---
$ws.Range("C11").Select()
